$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-23 from 45188 to 45189
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
